$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.668.92'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.973.98'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.11%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.969.34'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.488'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.147'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.50%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.08'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.437'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000217'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.05'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.448.62'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.34%  '
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.693.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.964.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '455.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.663'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  +1.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.83'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '54.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.72'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '450.29'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.166.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0771'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0376'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.115'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.99'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.39'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.97%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.241'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.28%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '117.53'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0490'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -10.00%  '
$ws.Range('E51').Value = '  +9.70%  '
